$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 417.3096593333333
$ws.Range("N2").Value = 1251.928978
$ws.Range("O2").Value = 0.8277575129472603
$ws.Range("P2").Value = 0.8277575129472603
$ws.Range("Q2").Value = 51017.36404212288
$ws.Range("R2").Value = 459156.276379106
$ws.Range("S2").Value = 0.1646111174365852
$ws.Range("T2").Value = 0.1646111174365852

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.194815333333333
$ws.Range("N3").Value = 3.584446
$ws.Range("O3").Value = 0.002369984366839822
$ws.Range("P3").Value = 0.002369984366839822
$ws.Range("Q3").Value = 146.06977686823
$ws.Range("R3").Value = 1314.62799181407
$ws.Range("S3").Value = 0.0004713044204741606
$ws.Range("T3").Value = 0.0004713044204741605

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 58.96764933333333
$ws.Range("N4").Value = 176.902948
$ws.Range("O4").Value = 0.1169656960121252
$ws.Range("P4").Value = 0.1169656960121252
$ws.Range("Q4").Value = 7208.97291846274
$ws.Range("R4").Value = 64880.75626616467
$ws.Range("S4").Value = 0.02326025873658316
$ws.Range("T4").Value = 0.02326025873658316

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.67269233333333
$ws.Range("N5").Value = 80.01807700000001
$ws.Range("O5").Value = 0.05290680667377473
$ws.Range("P5").Value = 0.05290680667377473
$ws.Range("Q5").Value = 3260.817055917385
$ws.Range("R5").Value = 29347.35350325647
$ws.Range("S5").Value = 0.01052125583924036
$ws.Range("T5").Value = 0.01052125583924036

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 417.3096593333333
$ws.Range("N6").Value = 1251.928978
$ws.Range("O6").Value = 0.8277575129472603
$ws.Range("P6").Value = 0.8277575129472603
$ws.Range("Q6").Value = 55312.20933753451
$ws.Range("R6").Value = 497809.8840378107
$ws.Range("S6").Value = 0.1784687382009831
$ws.Range("T6").Value = 0.1784687382009831

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.194815333333333
$ws.Range("N7").Value = 3.584446
$ws.Range("O7").Value = 0.002369984366839822
$ws.Range("P7").Value = 0.002369984366839822
$ws.Range("Q7").Value = 158.3665135923455
$ws.Range("R7").Value = 1425.29862233111
$ws.Range("S7").Value = 0.0005109807073812786
$ws.Range("T7").Value = 0.0005109807073812785

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 58.96764933333333
$ws.Range("N8").Value = 176.902948
$ws.Range("O8").Value = 0.1169656960121252
$ws.Range("P8").Value = 0.1169656960121252
$ws.Range("Q8").Value = 7815.853026930242
$ws.Range("R8").Value = 70342.67724237218
$ws.Range("S8").Value = 0.02521840013962368
$ws.Range("T8").Value = 0.02521840013962368

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.67269233333333
$ws.Range("N9").Value = 80.01807700000001
$ws.Range("O9").Value = 0.05290680667377473
$ws.Range("P9").Value = 0.05290680667377473
$ws.Range("Q9").Value = 3535.325648329994
$ws.Range("R9").Value = 31817.93083496995
$ws.Range("S9").Value = 0.01140697714200454
$ws.Range("T9").Value = 0.01140697714200454

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 417.3096593333333
$ws.Range("N10").Value = 1251.928978
$ws.Range("O10").Value = 0.8277575129472603
$ws.Range("P10").Value = 0.8277575129472603
$ws.Range("Q10").Value = 133576.4070946107
$ws.Range("R10").Value = 1202187.663851496
$ws.Range("S10").Value = 0.4309936831870296
$ws.Range("T10").Value = 0.4309936831870296

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.194815333333333
$ws.Range("N11").Value = 3.584446
$ws.Range("O11").Value = 0.002369984366839822
$ws.Range("P11").Value = 0.002369984366839822
$ws.Range("Q11").Value = 382.4477478503168
$ws.Range("R11").Value = 3442.029730652851
$ws.Range("S11").Value = 0.00123399458824973
$ws.Range("T11").Value = 0.00123399458824973

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 58.96764933333333
$ws.Range("N12").Value = 176.902948
$ws.Range("O12").Value = 0.1169656960121252
$ws.Range("P12").Value = 0.1169656960121252
$ws.Range("Q12").Value = 18874.92071318182
$ws.Range("R12").Value = 169874.2864186364
$ws.Range("S12").Value = 0.0609012607464092
$ws.Range("T12").Value = 0.0609012607464092

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 26.67269233333333
$ws.Range("N13").Value = 80.01807700000001
$ws.Range("O13").Value = 0.05290680667377473
$ws.Range("P13").Value = 0.05290680667377473
$ws.Range("Q13").Value = 8537.646636596908
$ws.Range("R13").Value = 76838.81972937218
$ws.Range("S13").Value = 0.02754731804584313
$ws.Range("T13").Value = 0.02754731804584313

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 417.3096593333333
$ws.Range("N14").Value = 1251.928978
$ws.Range("O14").Value = 0.8277575129472603
$ws.Range("P14").Value = 0.8277575129472603
$ws.Range("Q14").Value = 16638.09160458967
$ws.Range("R14").Value = 149742.824441307
$ws.Range("S14").Value = 0.0536839741226624
$ws.Range("T14").Value = 0.0536839741226624

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.194815333333333
$ws.Range("N15").Value = 3.584446
$ws.Range("O15").Value = 0.002369984366839822
$ws.Range("P15").Value = 0.002369984366839822
$ws.Range("Q15").Value = 47.63715989303111
$ws.Range("R15").Value = 428.73443903728
$ws.Range("S15").Value = 0.0001537046507346527
$ws.Range("T15").Value = 0.0001537046507346527

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 58.96764933333333
$ws.Range("N16").Value = 176.902948
$ws.Range("O16").Value = 0.1169656960121252
$ws.Range("P16").Value = 0.1169656960121252
$ws.Range("Q16").Value = 2351.033889037404
$ws.Range("R16").Value = 21159.30500133664
$ws.Range("S16").Value = 0.007585776389509129
$ws.Range("T16").Value = 0.007585776389509129

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 26.67269233333333
$ws.Range("N17").Value = 80.01807700000001
$ws.Range("O17").Value = 0.05290680667377473
$ws.Range("P17").Value = 0.05290680667377473
$ws.Range("Q17").Value = 1063.437398242818
$ws.Range("R17").Value = 9570.93658418536
$ws.Range("S17").Value = 0.003431255646686698
$ws.Range("T17").Value = 0.003431255646686698
